$d = $word.ActiveDocument

# Remove the first two empty "Title"-styled paragraphs that precede the
# third (kept) empty "Title" paragraph / table.
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()
